$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
$ws.Range("N6").Value = 8
$ws.Range("G12").Value = 2.35
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 2.8
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 9
$ws.Range("Z12").Value = 23
$ws.Range("AA12").Value = 21
$ws.Range("AC12").Value = 9
$ws.Range("AN12").Value = 4.5
$ws.Range("AW12").Value = 4.75
$ws.Range("AM15").Value = 900
$ws.Range("Q23").Value = 1.53
$ws.Range("R23").Value = 2.4
$ws.Range("G27").Value = 2.12
$ws.Range("H27").Value = 2.85
$ws.Range("I27").Value = 3.75
$ws.Range("J27").Value = 2.8
$ws.Range("K27").Value = 1.88
$ws.Range("L27").Value = 4.3
$ws.Range("N27").Value = 6
$ws.Range("P27").Value = 2.32
$ws.Range("Q27").Value = 2.37
$ws.Range("R27").Value = 1.45
$ws.Range("U27").Value = 2.02
$ws.Range("W27").Value = 5.6
$ws.Range("X27").Value = 9
$ws.Range("Y27").Value = 9.25
$ws.Range("Z27").Value = 20
$ws.Range("AA27").Value = 21
$ws.Range("AC27").Value = 6.2
$ws.Range("AD27").Value = 5.7
$ws.Range("AE27").Value = 17.5
$ws.Range("AF27").Value = 110
$ws.Range("AG27").Value = 8.25
$ws.Range("AH27").Value = 19
$ws.Range("AI27").Value = 13
$ws.Range("AJ27").Value = 60
$ws.Range("AL27").Value = 55
$ws.Range("AN27").Value = 3.8
$ws.Range("AO27").Value = 11.5
$ws.Range("AP27").Value = 23
$ws.Range("AQ27").Value = 50
$ws.Range("AS27").Value = 350
$ws.Range("AU27").Value = 7.5
$ws.Range("AV27").Value = 80
$ws.Range("AW27").Value = 5.3
$ws.Range("AX27").Value = 22
$ws.Range("AZ27").Value = 120
